$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add three new release rows to the "releases" sheet:
#   - "September"          (promotional release, no album) -> 2018-04-13
#   - "Lover (Remix)"       (promotional release, no album) -> 2019-11-13
#   - "Three Sad Virgins"   (promotional release, no album) -> 2021-11-13
# Each is inserted as a new row so the rows that follow shift down, matching
# the target layout exactly.
# ---------------------------------------------------------------------------

# 1) "September" belongs right before the existing "Babe" row (old row 73).
$ws.Rows.Item(73).Insert()
$ws.Range("B73").Value = "September"
$ws.Range("D73").Value = 43203

# 2) "Lover (Remix)" belongs right before the existing "Beautiful Ghosts" row
#    (old row 80, now row 81 after the insert above).
$ws.Rows.Item(81).Insert()
$ws.Range("B81").Value = "Lover (Remix)"
$ws.Range("D81").Value = 43782

# 3) "Three Sad Virgins" belongs right before the existing
#    "The Joker And The Queen" row (old row 104, now row 106 after the two
#    inserts above).
$ws.Rows.Item(106).Insert()
$ws.Range("B106").Value = "Three Sad Virgins"
$ws.Range("D106").Value = 44513

# ---------------------------------------------------------------------------
# The worksheet's AutoFilter range needs to grow from A1:D99 to A1:D101 (it
# only absorbs the two inserts that landed inside/above its old bound; the
# third insert, further down the sheet, falls outside it). Re-applying
# AutoFilter on a range always snaps to the full contiguous used range of the
# addressed columns, so temporarily clear the rows beyond the target bound,
# reapply the filter, then restore their contents.
# ---------------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row()
$saved = @()
for ($r = 102; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()
    $saved += ,@($a, $b, $c, $d)
}

$ws.Range("A102:D" + $lastRow).ClearContents() | Out-Null
$ws.AutoFilterMode = $false
$ws.Range("A1:D101").AutoFilter() | Out-Null

for ($i = 0; $i -lt $saved.Count; $i++) {
    $r = 102 + $i
    $row = $saved[$i]
    if ($row[0]) { $ws.Cells.Item($r, 1).Value = $row[0] }
    if ($row[1]) { $ws.Cells.Item($r, 2).Value = $row[1] }
    if ($row[2]) { $ws.Cells.Item($r, 3).Value = $row[2] }
    if ($row[3]) { $ws.Cells.Item($r, 4).Value = $row[3] }
}

# Keep the workbook-level _FilterDatabase defined name in sync with the new
# AutoFilter bound.
foreach ($n in $wb.Names) {
    if ($n.Name() -eq "releases!_FilterDatabase") {
        $n.RefersTo = "=releases!`$A`$1:`$D`$101"
    }
}

# Match the final active selection from the authored edit.
$ws.Range("B81").Select() | Out-Null
